# Update the "取得日時" (acquired timestamp) column for the existing
# lancers job rows from 2025-11-18 01:20:10 to 2025-11-18 01:49:18.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-11-18 01:20:10"
$newTimestamp = "2025-11-18 01:49:18"

for ($row = 2; $row -le 12; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    if ($cell.Text -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
